$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (values chosen so the resulting quantized/stored
# width matches the target stored width as closely as the COM pixel-grid
# rounding allows: target 2.140625 -> stored 2.166667, 3.140625 -> 3.166667,
# 5.7109375 -> 5.666667)
$ws.Columns.Item(4).ColumnWidth = 1.25
$ws.Columns.Item(7).ColumnWidth = 2.25
$ws.Columns.Item(9).ColumnWidth = 4.8
$ws.Columns.Item(10).ColumnWidth = 4.8

# Update row 1 values
$ws.Range("B1").Value = 3
$ws.Range("C1").Value = 33
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 31
$ws.Range("F1").Value = 13
$ws.Range("G1").Value = 11
$ws.Range("H1").Value = 21
$ws.Range("I1").Value = 0.077
$ws.Range("J1").Value = 0.011
$ws.Range("K1").Value = 0.069
